# Updated symbol list on Thu Feb  2 07:24:32 UTC 2023 with GitHub Actions
#
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) figures for
# the crypto-exchange-token rows on Sheet1. Values are written as literal
# text (matching the sheet's existing inlineStr-style cells, e.g. "6.99%"
# rather than the number 0.0699) by forcing a Text number format for the
# write and then restoring the cell to the built-in "Normal" style so no
# other formatting changes are introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($address, $value) {
    $cell = $ws.Range($address)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue "D2" "329.69"
Set-TextValue "E2" "6.99%"
Set-TextValue "D3" "39.96"
Set-TextValue "E3" "6.50%"
Set-TextValue "D4" "5.271"
Set-TextValue "E4" "2.46%"
Set-TextValue "D5" "0.08088"
Set-TextValue "E5" "2.98%"
Set-TextValue "D6" "4.509"
Set-TextValue "E6" "1.76%"
Set-TextValue "D7" "8.646"
Set-TextValue "E7" "5.03%"
Set-TextValue "D8" "1.930"
Set-TextValue "E8" "1.59%"
Set-TextValue "D10" "0.9373"
Set-TextValue "E10" "0.22%"
Set-TextValue "E11" "21.66%"
Set-TextValue "D12" "0.1979"
Set-TextValue "E12" "1.99%"
Set-TextValue "D13" "0.09090"
Set-TextValue "E13" "1.07%"
Set-TextValue "D14" "0.03503"
Set-TextValue "E14" "4.77%"
Set-TextValue "D15" "0.09590"
Set-TextValue "E15" "-0.01%"
Set-TextValue "D16" "0.001391"
Set-TextValue "E16" "0.99%"
Set-TextValue "D17" "0.006525"
Set-TextValue "E17" "5.94%"
Set-TextValue "D18" "3.366"
Set-TextValue "E18" "-6.73%"
Set-TextValue "D19" "0.3520"
Set-TextValue "E19" "3.21%"
Set-TextValue "D20" "6.768"
Set-TextValue "E20" "5.56%"
Set-TextValue "D21" "0.1313"
Set-TextValue "E21" "2.51%"
Set-TextValue "E22" "10.92%"
Set-TextValue "D23" "0.04429"
Set-TextValue "E23" "0.97%"
Set-TextValue "D24" "0.001222"
Set-TextValue "E24" "-1.03%"
Set-TextValue "D25" "0.004305"
Set-TextValue "E25" "-5.72%"
Set-TextValue "D26" "0.0001291"
Set-TextValue "E26" "-0.83%"
Set-TextValue "D27" "0.0003992"
Set-TextValue "E27" "-0.01%"
Set-TextValue "D39" "0.02495"
Set-TextValue "E39" "12.21%"
Set-TextValue "D40" "0.05224"
Set-TextValue "E40" "3.48%"
Set-TextValue "D41" "0.007725"
Set-TextValue "E41" "3.77%"
Set-TextValue "D42" "0.1428"
Set-TextValue "E42" "5.79%"
Set-TextValue "D43" "0.009230"
Set-TextValue "E43" "5.60%"
Set-TextValue "D44" "0.002171"
Set-TextValue "E44" "1.82%"
Set-TextValue "D45" "0.009320"
Set-TextValue "E45" "14.81%"
Set-TextValue "D46" "0.00006649"
Set-TextValue "E46" "1.55%"
Set-TextValue "D47" "0.00000000750"
Set-TextValue "D48" "0.003331"
Set-TextValue "E48" "16.45%"
Set-TextValue "E49" "148.02%"
Set-TextValue "D50" "0.00002101"
Set-TextValue "D51" "0.0002001"
